$wb = $excel.ActiveWorkbook

# Work on the "Repayment schedule" tab (4th sheet): insert a new blank
# column before the existing "N" column (shifting Late/heading/Outstanding
# one column to the right), set its width, move the selection, and make
# this sheet the active one (matches activeTab on the workbook + the
# tabSelected flag moving from "NewLoanInput" to this sheet).
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = 10.166666666666666

$ws.Range("L16").Select() | Out-Null
